# Generate Report for Handoff
# Refresh the "Latest Handoff" timestamps for rows that are queued for
# handoff (and the one failed-transform row) across the Overview sheet
# and each per-language detail sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

foreach ($r in $rows) {
    $overview.Range("D$r").Value = "2016-03-21 12:25:00"
    $zhcn.Range("E$r").Value = "2016-03-21 12:24:56"
    $dede.Range("E$r").Value = "2016-03-21 12:25:00"
}
